$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "76×76=5776" "24×81=1944"
Replace-Text "74×48=3552" "75×65=4875"
Replace-Text "79×12=948" "50×17=850"
Replace-Text "36×39=1404" "52×33=1716"
Replace-Text "21×44=924" "49×14=686"
Replace-Text "41×81=3321" "90×79=7110"
Replace-Text "39×33=1287" "26×28=728"
Replace-Text "20×22=440" "87×17=1479"
Replace-Text "33×73=2409" "49×36=1764"
Replace-Text "46×23=1058" "90×28=2520"
Replace-Text "91×32=2912" "76×91=6916"
Replace-Text "51×35=1785" "49×55=2695"
Replace-Text "53×76=4028" "75×42=3150"
Replace-Text "25×33=825" "57×45=2565"
Replace-Text "39×36=1404" "49×67=3283"
Replace-Text "61×13=793" "81×92=7452"
Replace-Text "27×65=1755" "53×46=2438"
Replace-Text "67×12=804" "96×56=5376"
Replace-Text "18×97=1746" "53×34=1802"
Replace-Text "87×55=4785" "80×72=5760"
Replace-Text "94×91=8554" "53×62=3286"
Replace-Text "34×98=3332" "21×39=819"
Replace-Text "27×93=2511" "44×71=3124"
Replace-Text "17×27=459" "85×45=3825"
Replace-Text "83×42=3486" "80×22=1760"
